# Error correction (Jeuk <> Puurs distance was wrong).
#
# The "Mere en Jeuk" sheet's Puurs<->Jeuk travel distance (cell P9, the
# Jeuk column on the Puurs row of the CityDistance matrix) was wrong.
# Fixing it changes the solver's "best" route found on that sheet (column
# AA, the order in which cities are visited); AB/AC and the AB26 total
# recompute automatically from the AA/P9 inputs via their existing
# formulas. Some solver-engine scratch settings for that sheet are also
# tightened to match the new run.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Mere en Jeuk")

# --- The actual data bug fix: Puurs (row 9) <-> Jeuk (column P) distance ---
$ws.Range("P9").Value = 97.2

# --- New best tour order found by the solver after the fix (column AA) ---
$tour = @{
  2  = 21
  3  = 7
  4  = 19
  5  = 14.999999999999998
  6  = 1
  7  = 5
  8  = 12
  9  = 4
  10 = 14
  11 = 13
  12 = 24
  13 = 3
  14 = 17
  15 = 9
  16 = 11
  17 = 20
  18 = 23
  19 = 8
  20 = 18
  21 = 16
  22 = 10
  23 = 22
  24 = 6
  25 = 2
}

foreach ($row in $tour.Keys) {
    $ws.Range("AA$row").Value = $tour[$row]
}

# --- Solver scratch settings (localSheetId=2, this sheet) tightened to match ---
$names = $wb.Names
for ($i = 1; $i -le $names.Count; $i++) {
    $item = $names.Item($i)
    switch ($item.Name) {
        "Mere en Jeuk!solver_mip" { $item.RefersTo = "=9999999" }
        "Mere en Jeuk!solver_mni" { $item.RefersTo = "=300" }
        "Mere en Jeuk!solver_nod" { $item.RefersTo = "=9999999" }
        "Mere en Jeuk!solver_pre" { $item.RefersTo = "=0.01" }
        "Mere en Jeuk!solver_tim" { $item.RefersTo = "=600" }
    }
}

# --- View state: activate the sheet and select the cells the fix highlighted ---
$ws.Activate()
$ws.Range("AC15:AC18").Select()
